# Apply cryptos list update (prices & % changes) per commit:
# 'Updated cryptos list on Wed Nov 29 05:54:54 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.010.70'
$ws.Range('D3').Value = '2.052.59'
$ws.Range('E3').Value = '  +2.58%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.42'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.52'
$ws.Range('E7').Value = '  +7.64%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +3.55%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0811'
$ws.Range('E10').Value = '  +4.64%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').Value = '2.354.36'
$ws.Range('E12').Value = '  +2.39%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.59'
$ws.Range('E13').Value = '  +4.95%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.83'
$ws.Range('E14').Value = '  +5.74%  '
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '2.042.92'
$ws.Range('E17').Value = '  +2.90%  '
$ws.Range('D18').Value = '37.909.92'
$ws.Range('E18').Value = '  +2.83%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.34'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.67'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  +3.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '223.94'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  +5.25%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '166.30'
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +3.47%  '
$ws.Range('E28').Value = '  +6.47%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.96'
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('E31').Value = '  +2.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  +11.10%  '
$ws.Range('E34').Value = '  +3.48%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0607'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.32'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.06'
$ws.Range('E37').Value = '  +14.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.27'
$ws.Range('E38').Value = '  +5.95%  '
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '1.518.88'
$ws.Range('E40').Value = '  +4.72%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '97.43'
$ws.Range('E41').Value = '  +3.40%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  +4.84%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0217'
$ws.Range('E43').Value = '  +2.94%  '
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0921'
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.06'
$ws.Range('E47').Value = '  +15.65%  '
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.97'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '2.242.11'
$ws.Range('E51').Value = '  +2.26%  '
